$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Variável" column before the existing B ("Valor" shifts to C),
# then insert a new "Colocação" column after it (becomes D).
$ws.Columns("B").Insert()
$ws.Columns("D").Insert()

# Row 1 headers
$ws.Range("B1").Value = "Variável"
$ws.Range("C1").Value = "Valor"
$ws.Range("D1").Value = "Colocação"

# Data rows: column B gets the repeated "Diferença 2021-2012" label,
# column D gets the ranking (only for rows 2-8).
$labels = "Diferença 2021-2012"

$ws.Range("B2").Value = $labels
$ws.Range("B3").Value = $labels
$ws.Range("B4").Value = $labels
$ws.Range("B5").Value = $labels
$ws.Range("B6").Value = $labels
$ws.Range("B7").Value = $labels
$ws.Range("B8").Value = $labels
$ws.Range("B9").Value = $labels
$ws.Range("B10").Value = $labels

$ws.Range("D2").Value = "1º"
$ws.Range("D3").Value = "2º"
$ws.Range("D4").Value = "3º"
$ws.Range("D5").Value = "4º"
$ws.Range("D6").Value = "5º"
$ws.Range("D7").Value = "6º"
$ws.Range("D8").Value = "12º"
